$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H125").Value = 3136.2727
$ws_ALC.Range("I125").Value = 1099
$ws_ALC.Range("J125").Value = 3340
$ws_ALC.Range("K125").Value = 9891
$ws_ALC.Range("L125").Value = 30060
$ws_ALC.Range("M125").Value = -7431
$ws_ALC.Range("N125").Value = -34980
$ws_ALC.Range("H132").Value = 14640841
$ws_ALC.Range("I132").Value = 18871166
$ws_ALC.Range("J132").Value = 627889.3
$ws_ALC.Range("K132").Value = 56613498
$ws_ALC.Range("L132").Value = 1883667.9
$ws_ALC.Range("M132").Value = -56610968
$ws_ALC.Range("N132").Value = -1888727.9
$ws_ALC.Range("H137").Value = 2310.093
$ws_ALC.Range("I137").Value = 1150.1562
$ws_ALC.Range("J137").Value = 5684.4546
$ws_ALC.Range("K137").Value = 3450.4686
$ws_ALC.Range("L137").Value = 17053.3638
$ws_ALC.Range("M137").Value = -900.4685999999997
$ws_ALC.Range("N137").Value = -22153.3638
$ws_ALC.Range("H138").Value = 2497.62
$ws_ALC.Range("I138").Value = 860
$ws_ALC.Range("J138").Value = 3268.2646
$ws_ALC.Range("K138").Value = 2580
$ws_ALC.Range("L138").Value = 9804.793799999999
$ws_ALC.Range("M138").Value = 2560
$ws_ALC.Range("N138").Value = -20084.7938

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H15").Value = 0
$ws_ARM.Range("J15").Value = 0
$ws_ARM.Range("L15").ClearContents()
$ws_ARM.Range("N15").Value = 0
$ws_ARM.Range("H61").Value = 1077.6129
$ws_ARM.Range("I61").Value = 928.7368
$ws_ARM.Range("J61").Value = 1313.3334
$ws_ARM.Range("K61").Value = 928.7368
$ws_ARM.Range("L61").Value = 1313.3334
$ws_ARM.Range("M61").Value = -716.7368
$ws_ARM.Range("N61").Value = -1737.3334
$ws_ARM.Range("H74").Value = 2300.0178
$ws_ARM.Range("I74").Value = 2246.5334
$ws_ARM.Range("K74").Value = 2246.5334
$ws_ARM.Range("M74").Value = -1372.5334
$ws_ARM.Range("H77").Value = 2300.0178
$ws_ARM.Range("I77").Value = 2246.5334
$ws_ARM.Range("K77").Value = 11232.667
$ws_ARM.Range("M77").Value = -6864.666999999999
$ws_ARM.Range("H80").Value = 26586.4
$ws_ARM.Range("J80").Value = 26586.4
$ws_ARM.Range("L80").Value = 26586.4
$ws_ARM.Range("N80").Value = -28582.4
$ws_ARM.Range("H83").Value = 26586.4
$ws_ARM.Range("J83").Value = 26586.4
$ws_ARM.Range("L83").Value = 79759.20000000001
$ws_ARM.Range("N83").Value = -89743.20000000001
$ws_ARM.Range("H132").Value = 2055.8044
$ws_ARM.Range("I132").Value = 1351.5588
$ws_ARM.Range("K132").Value = 4054.6764
$ws_ARM.Range("M132").Value = -1524.6764
$ws_ARM.Range("H136").Value = 1077.6129
$ws_ARM.Range("I136").Value = 928.7368
$ws_ARM.Range("J136").Value = 1313.3334
$ws_ARM.Range("K136").Value = 2786.2104
$ws_ARM.Range("L136").Value = 3940.0002
$ws_ARM.Range("M136").Value = -236.2103999999999
$ws_ARM.Range("N136").Value = -9040.0002

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H30").Value = 20000
$ws_BSM.Range("J30").Value = 20000
$ws_BSM.Range("L30").Value = 20000
$ws_BSM.Range("N30").Value = -20250
$ws_BSM.Range("H80").Value = 148.83333
$ws_BSM.Range("I80").Value = 83.75
$ws_BSM.Range("J80").Value = 192.22223
$ws_BSM.Range("K80").Value = 83.75
$ws_BSM.Range("L80").Value = 192.22223
$ws_BSM.Range("M80").Value = 914.25
$ws_BSM.Range("N80").Value = -2188.22223
$ws_BSM.Range("H83").Value = 148.83333
$ws_BSM.Range("I83").Value = 83.75
$ws_BSM.Range("J83").Value = 192.22223
$ws_BSM.Range("K83").Value = 418.75
$ws_BSM.Range("L83").Value = 961.11115
$ws_BSM.Range("M83").Value = 4573.25
$ws_BSM.Range("N83").Value = -10945.11115
$ws_BSM.Range("H99").Value = 3414.4443
$ws_BSM.Range("I99").Value = 1306.6666
$ws_BSM.Range("J99").Value = 4468.3335
$ws_BSM.Range("K99").Value = 1306.6666
$ws_BSM.Range("L99").Value = 4468.3335
$ws_BSM.Range("M99").Value = 191.3334
$ws_BSM.Range("N99").Value = -7464.3335
$ws_BSM.Range("H105").Value = 1765.5151
$ws_BSM.Range("I105").Value = 1634.84
$ws_BSM.Range("K105").Value = 1634.84
$ws_BSM.Range("M105").Value = 112.1600000000001
$ws_BSM.Range("H134").Value = 2566.4883
$ws_BSM.Range("I134").Value = 1404.8077
$ws_BSM.Range("J134").Value = 4343.1763
$ws_BSM.Range("K134").Value = 4214.4231
$ws_BSM.Range("L134").Value = 13029.5289
$ws_BSM.Range("M134").Value = -1679.4231
$ws_BSM.Range("N134").Value = -18099.5289

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 8066519.5
$ws_CRP.Range("I31").Value = 1104.0889
$ws_CRP.Range("K31").Value = 1104.0889
$ws_CRP.Range("M31").Value = -809.0889
$ws_CRP.Range("H34").Value = 8066519.5
$ws_CRP.Range("I34").Value = 1104.0889
$ws_CRP.Range("K34").Value = 1104.0889
$ws_CRP.Range("M34").Value = -902.0889
$ws_CRP.Range("H58").Value = 1831.0541
$ws_CRP.Range("I58").Value = 1581.5151
$ws_CRP.Range("K58").Value = 1581.5151
$ws_CRP.Range("M58").Value = -1378.5151
$ws_CRP.Range("H107").Value = 955.63635
$ws_CRP.Range("I107").Value = 591.5
$ws_CRP.Range("K107").Value = 591.5
$ws_CRP.Range("M107").Value = 1328.5
$ws_CRP.Range("H115").Value = 29460
$ws_CRP.Range("J115").Value = 29460
$ws_CRP.Range("L115").Value = 29460
$ws_CRP.Range("N115").Value = -31810
$ws_CRP.Range("H132").Value = 3509.6155
$ws_CRP.Range("I132").Value = 3540.6
$ws_CRP.Range("J132").Value = 3467.3635
$ws_CRP.Range("K132").Value = 10621.8
$ws_CRP.Range("L132").Value = 10402.0905
$ws_CRP.Range("M132").Value = -8091.799999999999
$ws_CRP.Range("N132").Value = -15462.0905
$ws_CRP.Range("H134").Value = 3042.1355
$ws_CRP.Range("I134").Value = 3157.0488
$ws_CRP.Range("J134").Value = 2780.389
$ws_CRP.Range("K134").Value = 9471.1464
$ws_CRP.Range("L134").Value = 8341.167000000001
$ws_CRP.Range("M134").Value = -6936.1464
$ws_CRP.Range("N134").Value = -13411.167
$ws_CRP.Range("H136").Value = 1831.0541
$ws_CRP.Range("I136").Value = 1581.5151
$ws_CRP.Range("K136").Value = 4744.5453
$ws_CRP.Range("M136").Value = -2194.5453

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H5").Value = 1346.359
$ws_CUL.Range("I5").Value = 466.3889
$ws_CUL.Range("K5").Value = 1399.1667
$ws_CUL.Range("M5").Value = -1287.1667
$ws_CUL.Range("H80").Value = 7086.727
$ws_CUL.Range("J80").Value = 6995.4
$ws_CUL.Range("L80").Value = 20986.2
$ws_CUL.Range("N80").Value = -22858.2
$ws_CUL.Range("H83").Value = 7086.727
$ws_CUL.Range("J83").Value = 6995.4
$ws_CUL.Range("L83").Value = 62958.6
$ws_CUL.Range("N83").Value = -72318.60000000001
$ws_CUL.Range("H110").Value = 3999
$ws_CUL.Range("I110").Value = 3999
$ws_CUL.Range("K110").Value = 11997
$ws_CUL.Range("M110").Value = -7907
$ws_CUL.Range("H122").Value = 3153.7188
$ws_CUL.Range("I122").Value = 602.625
$ws_CUL.Range("J122").Value = 4004.0833
$ws_CUL.Range("K122").Value = 5423.625
$ws_CUL.Range("L122").Value = 36036.7497
$ws_CUL.Range("M122").Value = -2973.625
$ws_CUL.Range("N122").Value = -40936.7497
$ws_CUL.Range("H135").Value = 1346.359
$ws_CUL.Range("I135").Value = 466.3889
$ws_CUL.Range("K135").Value = 4197.5001
$ws_CUL.Range("M135").Value = -1662.5001
$ws_CUL.Range("H137").Value = 2518.7026
$ws_CUL.Range("I137").Value = 678.9091
$ws_CUL.Range("J137").Value = 3297.077
$ws_CUL.Range("K137").Value = 2036.7273
$ws_CUL.Range("L137").Value = 9891.231
$ws_CUL.Range("M137").Value = 3063.2727
$ws_CUL.Range("N137").Value = -20091.231

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H132").Value = 2294.195
$ws_GSM.Range("I132").Value = 1394.4286
$ws_GSM.Range("J132").Value = 4232.154
$ws_GSM.Range("K132").Value = 4183.2858
$ws_GSM.Range("L132").Value = 12696.462
$ws_GSM.Range("M132").Value = -1653.2858
$ws_GSM.Range("N132").Value = -17756.462

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H7").Value = 4239.3887
$ws_LTW.Range("I7").Value = 3330
$ws_LTW.Range("J7").Value = 5376.125
$ws_LTW.Range("K7").Value = 3330
$ws_LTW.Range("L7").Value = 5376.125
$ws_LTW.Range("M7").Value = -3218
$ws_LTW.Range("N7").Value = -5600.125
$ws_LTW.Range("H126").Value = 4239.3887
$ws_LTW.Range("I126").Value = 3330
$ws_LTW.Range("J126").Value = 5376.125
$ws_LTW.Range("K126").Value = 9990
$ws_LTW.Range("L126").Value = 16128.375
$ws_LTW.Range("M126").Value = -7520
$ws_LTW.Range("N126").Value = -21068.375
$ws_LTW.Range("H132").Value = 9358.604499999999
$ws_LTW.Range("I132").Value = 9375.532999999999
$ws_LTW.Range("J132").Value = 9319.538
$ws_LTW.Range("K132").Value = 28126.599
$ws_LTW.Range("L132").Value = 27958.614
$ws_LTW.Range("M132").Value = -25596.599
$ws_LTW.Range("N132").Value = -33018.614

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H132").Value = 7753882
$ws_WVR.Range("I132").Value = 1423.742
$ws_WVR.Range("J132").Value = 27781066
$ws_WVR.Range("K132").Value = 4271.226
$ws_WVR.Range("L132").Value = 83343198
$ws_WVR.Range("M132").Value = -1741.226
$ws_WVR.Range("N132").Value = -83348258
